$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force literal text (no date/number auto-conversion) for columns that
# look like dates or plain numbers, then strip the resulting number-format
# style so the cell ends up unstyled (matching the rest of the data rows).
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "2023-06-13"
$ws.Range("A44").ClearFormats()

$ws.Range("B44").Value = "22:23:37"
$ws.Range("C44").Value = "Tuesday"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = 121136
$ws.Range("F44").Value = 135066
$ws.Range("G44").Value = 161256
$ws.Range("H44").Value = 132418
$ws.Range("I44").Value = 176376
$ws.Range("J44").Value = 113946
$ws.Range("K44").Value = 202388
$ws.Range("L44").Value = 223251
$ws.Range("M44").Value = 173942
$ws.Range("N44").Value = 101590
$ws.Range("O44").Value = 38916
$ws.Range("P44").Value = 34068
$ws.Range("Q44").Value = 51460
$ws.Range("R44").Value = -1
$ws.Range("S44").Value = 36980
$ws.Range("T44").Value = -1
